$d = $word.ActiveDocument

# List of (old, new) text replacements, in document order.
$replacements = @(
    @("2025-07-07 Monday", "2025-07-08 Tuesday"),
    @("16÷8=", "57÷4="),
    @("75÷8=", "95÷3="),
    @("41÷7=", "69÷8="),
    @("25÷2=", "60÷7="),
    @("34÷2=", "70÷2="),
    @("53÷9=", "59÷3="),
    @("67÷8=", "93÷8="),
    @("46÷7=", "92÷9="),
    @("47÷8=", "81÷6="),
    @("52÷6=", "50÷4="),
    @("80÷5=", "14÷6="),
    @("54÷4=", "20÷9="),
    @("57÷9=", "70÷7="),
    @("34÷9=", "16÷8="),
    @("70÷4=", "66÷9="),
    @("21÷8=", "85÷8="),
    @("44÷9=", "10÷4="),
    @("96÷7=", "45÷8="),
    @("55÷5=", "99÷8="),
    @("87÷3=", "45÷3="),
    @("19÷8=", "62÷6="),
    @("29÷5=", "67÷3="),
    @("26÷4=", "45÷6="),
    @("33÷8=", "30÷9="),
    @("52÷9=", "55÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
